# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Updated K values per game row (row 2 = most recent game, row 38 = oldest).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 2
    3  = 0
    4  = 3
    5  = 0
    6  = 2
    7  = 1
    8  = 0
    9  = 0
    10 = 1
    11 = 3
    12 = 0
    13 = 2
    14 = 5
    15 = 2
    16 = 5
    17 = 3
    18 = 1
    19 = 8
    20 = 3
    21 = 3
    22 = 5
    23 = 5
    24 = 2
    25 = 3
    26 = 6
    27 = 2
    28 = 7
    29 = 5
    30 = 2
    31 = 4
    32 = 5
    33 = 3
    34 = 3
    35 = 4
    36 = 3
    37 = 2
    38 = 4
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
